# Apply the 2023-04-22 cryptos data refresh (GitHub Actions scheduled update).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextCell($cellRef, $text) {
    # Leading apostrophe forces Excel to treat the value as literal text
    # (prevents auto-conversion of numeric-looking strings like "0.9999"
    # into real numbers), then reset the style so no quote-prefix format
    # lingers on the cell.
    $r = $ws.Range($cellRef)
    $r.Value = "'" + $text
    $r.Style = "Normal"
}

Set-TextCell "D2" "27.327.90"
Set-TextCell "E2" "  -4.33%  "
Set-TextCell "D3" "1.860.81"
Set-TextCell "E3" "  -5.33%  "
Set-TextCell "D4" "0.9999"
Set-TextCell "D5" "323.21"
Set-TextCell "E5" "  -0.18%  "
Set-TextCell "E6" "  -0.95%  "
Set-TextCell "D7" "0.4514"
Set-TextCell "E7" "  -6.54%  "
Set-TextCell "E8" "  -5.50%  "
Set-TextCell "D9" "48.12"
Set-TextCell "E9" "  -11.00%  "
Set-TextCell "D10" "0.07911"
Set-TextCell "E10" "  -7.29%  "
Set-TextCell "E11" "  -4.24%  "
Set-TextCell "D12" "21.41"
Set-TextCell "E12" "  -4.97%  "
Set-TextCell "D13" "1.873.53"
Set-TextCell "E13" "  -8.03%  "
Set-TextCell "B14" "Polkadot"
Set-TextCell "C14" "https://coinranking.com/coin/25W7FG7om+polkadot-dot"
Set-TextCell "D14" "5.882"
Set-TextCell "E14" "  -5.36%  "
Set-TextCell "B15" "Chainlink"
Set-TextCell "C15" "https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link"
Set-TextCell "D15" "7.157"
Set-TextCell "E15" "  -6.20%  "
Set-TextCell "D16" "0.9997"
Set-TextCell "E16" "  -1.25%  "
Set-TextCell "D17" "0.00001037"
Set-TextCell "E17" "  -3.66%  "
Set-TextCell "D18" "85.74"
Set-TextCell "E18" "  -6.19%  "
Set-TextCell "D19" "0.06523"
Set-TextCell "E19" "  -1.86%  "
Set-TextCell "D20" "17.07"
Set-TextCell "E20" "  -8.18%  "
Set-TextCell "D21" "1.000"
Set-TextCell "E21" "  -0.99%  "
Set-TextCell "D22" "5.516"
Set-TextCell "E22" "  -6.28%  "
Set-TextCell "D23" "27.333.86"
Set-TextCell "E23" "  -4.49%  "
Set-TextCell "D24" "10.81"
Set-TextCell "E24" "  -6.34%  "
Set-TextCell "D25" "2.269"
Set-TextCell "E25" "  -1.32%  "
Set-TextCell "D26" "2.084.92"
Set-TextCell "E26" "  -8.28%  "
Set-TextCell "D27" "152.30"
Set-TextCell "E27" "  -3.03%  "
Set-TextCell "D28" "19.75"
Set-TextCell "E28" "  -3.15%  "
Set-TextCell "D29" "2.065"
Set-TextCell "E29" "  -5.87%  "
Set-TextCell "D30" "5.500"
Set-TextCell "E30" "  -7.47%  "
Set-TextCell "D31" "120.58"
Set-TextCell "E31" "  -3.62%  "
Set-TextCell "D32" "1.491"
Set-TextCell "E32" "  +1.16%  "
Set-TextCell "D33" "0.09302"
Set-TextCell "E33" "  -4.48%  "
Set-TextCell "D34" "0.9373"
Set-TextCell "E34" "  -6.16%  "
Set-TextCell "E35" "  -2.77%  "
Set-TextCell "D36" "5.283"
Set-TextCell "E36" "  -6.95%  "
Set-TextCell "D37" "0.02232"
Set-TextCell "E37" "  -4.87%  "
Set-TextCell "D38" "0.05994"
Set-TextCell "E38" "  -4.58%  "
Set-TextCell "D39" "1.210"
Set-TextCell "E39" "  -4.05%  "
Set-TextCell "D40" "8.259"
Set-TextCell "E40" "  -10.34%  "
Set-TextCell "D41" "0.9997"
Set-TextCell "E41" "  -0.96%  "
Set-TextCell "D42" "0.5901"
Set-TextCell "E42" "  -5.90%  "
Set-TextCell "D43" "0.1888"
Set-TextCell "E43" "  -2.05%  "
Set-TextCell "D44" "10.15"
Set-TextCell "E44" "  -10.29%  "
Set-TextCell "D45" "1.269"
Set-TextCell "E45" "  -6.32%  "
Set-TextCell "D46" "0.5636"
Set-TextCell "E46" "  -5.97%  "
Set-TextCell "D47" "12.00"
Set-TextCell "E47" "  -8.92%  "
Set-TextCell "B48" "PancakeSwap"
Set-TextCell "C48" "https://coinranking.com/coin/ncYFcP709+pancakeswap-cake"
Set-TextCell "D48" "3.364"
Set-TextCell "E48" "  -1.55%  "
Set-TextCell "B49" "NEARProtocol"
Set-TextCell "C49" "https://coinranking.com/coin/DCrsaMv68+nearprotocol-near"
Set-TextCell "D49" "1.923"
Set-TextCell "E49" "  -7.72%  "
Set-TextCell "D50" "0.06799"
Set-TextCell "E50" "  -0.61%  "
Set-TextCell "D51" "108.15"
Set-TextCell "E51" "  -3.20%  "
